# "added login test cases (#1)"
#
# RUNMANAGER: rename the two existing test cases to their new, more
# descriptive names (the test descriptions / flags stay the same).
#
# DATA: the two test cases now each get two data rows (one per browser)
# with updated credentials, and the old extra "subscribe" row is removed.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("DATA")

# --- RUNMANAGER sheet: rename the test cases ---
$ws1.Range("A2").Value = "verifyThatAdminCanLogInWithValidCredentials"
$ws1.Range("A3").Value = "verifyThatAdminCannotLogInWithInvalidCredentials"

# --- DATA sheet: update the test data rows ---

# Row 2: valid-login test, chrome
$ws2.Range("A2").Value = "verifyThatAdminCanLogInWithValidCredentials"
$ws2.Range("F2").Value = "chrome"

# Row 3: valid-login test, firefox
$ws2.Range("A3").Value = "verifyThatAdminCanLogInWithValidCredentials"
$ws2.Range("E3").Value = "amuthan"
$ws2.Range("F3").Value = "firefox"

# Row 4: invalid-login test, chrome
$ws2.Range("A4").Value = "verifyThatAdminCannotLogInWithInvalidCredentials"
$ws2.Range("C4").Value = "admin12"
$ws2.Range("E4").Value = "sunil"
$ws2.Range("F4").Value = "chrome"

# Row 5: invalid-login test, firefox
$ws2.Range("A5").Value = "verifyThatAdminCannotLogInWithInvalidCredentials"
$ws2.Range("C5").Value = "admin12"
$ws2.Range("E5").Value = "sunil"
$ws2.Range("F5").Value = "firefox"

# Row 6 (old loginLogoutTest/subscribe row) no longer exists
$ws2.Rows.Item(6).Delete()
